$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (aa / b / a@b.com) becomes (s / t / s@t.com)
$ws.Range("A3").Value = "s"
$ws.Range("B3").Value = "t"
$ws.Range("C3").Value = "s@t.com"

# Remove the rows that are no longer present in the final table:
#  - original row 4 (ab / b / c@a.com)
#  - original rows 6 and 7 (aabbcc/user/aabbcc@gmail.com, abcabc/wick/abcabc@gmail.com)
# Delete bottom-up so the row numbers used below stay valid.
$ws.Rows(7).Delete()
$ws.Rows(6).Delete()
$ws.Rows(4).Delete()
